$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.518.66"
$ws.Range("E2").Value = "'  +5.62%  "
$ws.Range("D3").Value = "'2.056.50"
$ws.Range("E3").Value = "'  +4.16%  "
$ws.Range("E4").Value = "'  -0.12%  "
$ws.Range("D5").Value = "'252.54"
$ws.Range("E5").Value = "'  +3.22%  "
$ws.Range("D6").Value = "'0.649"
$ws.Range("E6").Value = "'  +2.17%  "
$ws.Range("D7").Value = "'66.42"
$ws.Range("E7").Value = "'  +16.45%  "
$ws.Range("E8").Value = "'  -0.07%  "
$ws.Range("E9").Value = "'  +6.76%  "
$ws.Range("D10").Value = "'59.82"
$ws.Range("E10").Value = "'  +2.05%  "
$ws.Range("D11").Value = "'0.0768"
$ws.Range("E11").Value = "'  +5.02%  "
$ws.Range("E12").Value = "'  +1.38%  "
$ws.Range("D13").Value = "'0.912"
$ws.Range("E13").Value = "'  -3.87%  "
$ws.Range("D14").Value = "'14.99"
$ws.Range("E14").Value = "'  +5.61%  "
$ws.Range("D15").Value = "'2.359.32"
$ws.Range("E15").Value = "'  +4.29%  "
$ws.Range("D16").Value = "'21.17"
$ws.Range("E16").Value = "'  +20.83%  "
$ws.Range("D17").Value = "'5.59"
$ws.Range("E17").Value = "'  +6.44%  "
$ws.Range("D18").Value = "'2.056.62"
$ws.Range("E18").Value = "'  +4.24%  "
$ws.Range("D19").Value = "'37.280.48"
$ws.Range("E19").Value = "'  +5.25%  "
$ws.Range("D20").Value = "'73.88"
$ws.Range("E20").Value = "'  +3.69%  "
$ws.Range("E21").Value = "'  +4.65%  "
$ws.Range("D22").Value = "'5.47"
$ws.Range("E22").Value = "'  +6.55%  "
$ws.Range("D23").Value = "'240.35"
$ws.Range("E23").Value = "'  +3.41%  "
$ws.Range("D24").Value = "'2.65"
$ws.Range("E24").Value = "'  +3.47%  "
$ws.Range("E25").Value = "'  -0.12%  "
$ws.Range("E26").Value = "'  +4.02%  "
$ws.Range("D27").Value = "'9.77"
$ws.Range("E27").Value = "'  +7.86%  "
$ws.Range("D28").Value = "'161.34"
$ws.Range("D29").Value = "'20.04"
$ws.Range("E29").Value = "'  +4.29%  "
$ws.Range("E30").Value = "'  +9.29%  "
$ws.Range("E31").Value = "'  +26.39%  "
$ws.Range("E32").Value = "'  +3.37%  "
$ws.Range("E33").Value = "'  +6.55%  "
$ws.Range("E34").Value = "'  +11.79%  "
$ws.Range("E35").Value = "'  +5.46%  "
$ws.Range("E36").Value = "'  +4.12%  "
$ws.Range("E37").Value = "'  -0.07%  "
$ws.Range("E38").Value = "'  +4.47%  "
$ws.Range("D39").Value = "'6.11"
$ws.Range("E39").Value = "'  +19.83%  "
$ws.Range("D40").Value = "'3.02"
$ws.Range("E40").Value = "'  +34.28%  "
$ws.Range("E41").Value = "'  +16.97%  "
$ws.Range("E42").Value = "'  +3.21%  "
$ws.Range("E43").Value = "'  +4.98%  "
$ws.Range("E44").Value = "'  +4.29%  "
$ws.Range("E45").Value = "'  +6.58%  "
$ws.Range("D46").Value = "'17.02"
$ws.Range("E46").Value = "'  +7.10%  "
$ws.Range("D47").Value = "'95.88"
$ws.Range("E47").Value = "'  +5.21%  "
$ws.Range("D48").Value = "'7.97"
$ws.Range("E48").Value = "'  +6.54%  "
$ws.Range("D49").Value = "'1.420.32"
$ws.Range("E49").Value = "'  +2.95%  "
$ws.Range("E50").Value = "'  +2.44%  "
$ws.Range("D51").Value = "'46.71"
$ws.Range("E51").Value = "'  +2.59%  "
